$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 274; this shifts the existing rows 274-361
# down to 275-362 (matching the dimension growing from A1:R361 to A1:R362).
$ws.Rows("274").Insert()

# Populate the newly inserted row 274 with the new price-report record.
$ws.Range("A274").Value = 9
$ws.Range("B274").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C274").Value = 'Metropolitana'
$ws.Range("D274").Value = 45215
$ws.Range("E274").Value = 13
$ws.Range("F274").Value = 100112003
$ws.Range("G274").Value = 'Ajo'
$ws.Range("H274").Value = 'Chino'
$ws.Range("I274").Value = 'Primera'
$ws.Range("J274").Value = 520
$ws.Range("K274").Value = 18000
$ws.Range("L274").Value = 20000
$ws.Range("M274").Value = 19000
$ws.Range("N274").Value = '$/caja 10 kilos'
$ws.Range("O274").Value = 'China'
$ws.Range("P274").Value = 1900
$ws.Range("Q274").Value = 10
$ws.Range("R274").Value = 'Hortaliza'
